$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Control"

# Header row values (A1:I1)
$ws.Range("A1").Value = "CUIT"
$ws.Range("B1").Value = "Controbuyente"
$ws.Range("C1").Value = "Fila"
$ws.Range("D1").Value = "RET 216"
$ws.Range("E1").Value = "RET 217"
$ws.Range("F1").Value = "RET 767"
$ws.Range("G1").Value = "SIRCREB"
$ws.Range("H1").Value = "SIFERE"
$ws.Range("I1").Value = "CABA - AGIP"

# Style the header row: white font on a dark blue (Accent1, Darker 50%) fill with thin black borders
$header = $ws.Range("A1:I1")
$header.Font.ThemeColor = 2
$header.Interior.ThemeColor = 5
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

# Column B width (best-fit for "Controbuyente")
$ws.Columns.Item(2).ColumnWidth = 13.592447916666666

# Selection shown when the file is opened
$ws.Range("D2").Select() | Out-Null
